$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1400
$ws.Range("I9").Value = 1625.25
$ws.Range("K9").Value = 1625.25
$ws.Range("M9").Value = -1456.25

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H46").Value = 5000
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 15000
$ws.Range("L46").ClearContents()
$ws.Range("N46").Value = 0
$ws.Range("M46").Value = -14881

$ws.Range("H60").Value = 5000
$ws.Range("I60").Value = 5000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 15000
$ws.Range("L60").ClearContents()
$ws.Range("N60").Value = 0
$ws.Range("M60").Value = -14516

$ws.Range("H62").Value = 5933.6665
$ws.Range("I62").Value = 5897.5
$ws.Range("K62").Value = 5897.5
$ws.Range("M62").Value = -5273.5

$ws.Range("H65").Value = 5933.6665
$ws.Range("I65").Value = 5897.5
$ws.Range("K65").Value = 29487.5
$ws.Range("M65").Value = -26367.5

$ws.Range("H135").Value = 1485.9231
$ws.Range("I135").Value = 1490.0834
$ws.Range("J135").Value = 1436
$ws.Range("K135").Value = 13410.7506
$ws.Range("L135").Value = 12924
$ws.Range("M135").Value = -10875.7506
$ws.Range("N135").Value = -17994

$ws.Range("H137").Value = 2296.889
$ws.Range("I137").Value = 1335.2
$ws.Range("J137").Value = 3499
$ws.Range("K137").Value = 4005.6
$ws.Range("L137").Value = 10497
$ws.Range("M137").Value = -1455.6
$ws.Range("N137").Value = -15597

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 532.5
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H74").Value = 1448
$ws.Range("I74").Value = 1548.5385
$ws.Range("K74").Value = 1548.5385
$ws.Range("M74").Value = -674.5385000000001

$ws.Range("H77").Value = 1448
$ws.Range("I77").Value = 1548.5385
$ws.Range("K77").Value = 7742.692500000001
$ws.Range("M77").Value = -3374.692500000001

$ws.Range("H102").Value = 1000
$ws.Range("J102").Value = 1000
$ws.Range("L102").Value = 1000
$ws.Range("N102").Value = -4244

$ws.Range("H132").Value = 1841.2368
$ws.Range("I132").Value = 1699.5151
$ws.Range("J132").Value = 2776.6
$ws.Range("K132").Value = 5098.5453
$ws.Range("L132").Value = 8329.799999999999
$ws.Range("M132").Value = -2568.5453
$ws.Range("N132").Value = -13389.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 532.5
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H20").Value = 3600
$ws.Range("I20").Value = 3600
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 3600
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0
$ws.Range("M20").Value = -3353

$ws.Range("H134").Value = 2382.5334
$ws.Range("I134").Value = 2382.5334
$ws.Range("K134").Value = 7147.600199999999
$ws.Range("M134").Value = -4612.600199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 38990
$ws.Range("I41").Value = 39000
$ws.Range("K41").Value = 39000
$ws.Range("M41").Value = -38572

$ws.Range("H50").Value = 46790
$ws.Range("J50").Value = 46790
$ws.Range("L50").Value = 46790
$ws.Range("N50").Value = -48040

$ws.Range("H60").Value = 36853.773
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 36853.773
$ws.Range("K60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("M60").Value = 36853.773
$ws.Range("N60").Value = -37875.773

$ws.Range("H68").Value = 34995
$ws.Range("J68").Value = 34995
$ws.Range("L68").Value = 34995
$ws.Range("N68").Value = -36493

$ws.Range("H71").Value = 34995
$ws.Range("J71").Value = 34995
$ws.Range("L71").Value = 104985
$ws.Range("N71").Value = -112473

$ws.Range("H74").Value = 34534.5
$ws.Range("J74").Value = 34534.5
$ws.Range("L74").Value = 34534.5
$ws.Range("N74").Value = -36282.5

$ws.Range("H77").Value = 34534.5
$ws.Range("J77").Value = 34534.5
$ws.Range("L77").Value = 103603.5
$ws.Range("N77").Value = -112339.5

$ws.Range("H134").Value = 4381.636
$ws.Range("I134").Value = 4381.636
$ws.Range("K134").Value = 13144.908
$ws.Range("M134").Value = -10609.908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 257.25
$ws.Range("I12").Value = 267.5
$ws.Range("J12").Value = 247
$ws.Range("K12").Value = 802.5
$ws.Range("L12").Value = 741
$ws.Range("M12").Value = -629.5
$ws.Range("N12").Value = -1087

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5545.1113
$ws.Range("I80").Value = 5001.3335
$ws.Range("J80").Value = 5817
$ws.Range("K80").Value = 5001.3335
$ws.Range("L80").Value = 5817
$ws.Range("M80").Value = -4003.3335
$ws.Range("N80").Value = -7813

$ws.Range("H83").Value = 5545.1113
$ws.Range("I83").Value = 5001.3335
$ws.Range("J83").Value = 5817
$ws.Range("K83").Value = 25006.6675
$ws.Range("L83").Value = 29085
$ws.Range("M83").Value = -20014.6675
$ws.Range("N83").Value = -39069

$ws.Range("H113").Value = 973.5
$ws.Range("I113").Value = 973.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 973.5
$ws.Range("L113").ClearContents()
$ws.Range("N113").Value = 0
$ws.Range("M113").Value = 1196.5

$ws.Range("H132").Value = 2687.4707
$ws.Range("I132").Value = 2198.2307
$ws.Range("K132").Value = 6594.6921
$ws.Range("M132").Value = -4064.6921

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3050
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 3866.6667
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 3866.6667
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -4242.6667

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H136").Value = 9999
$ws.Range("I136").Value = 9999
$ws.Range("K136").Value = 29997
$ws.Range("M136").Value = -27447

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10499.5
$ws.Range("I62").Value = 6499.5
$ws.Range("J62").Value = 14499.5
$ws.Range("K62").Value = 6499.5
$ws.Range("L62").Value = 14499.5
$ws.Range("M62").Value = -5875.5
$ws.Range("N62").Value = -15747.5

$ws.Range("H65").Value = 10499.5
$ws.Range("I65").Value = 6499.5
$ws.Range("J65").Value = 14499.5
$ws.Range("K65").Value = 32497.5
$ws.Range("L65").Value = 72497.5
$ws.Range("M65").Value = -29377.5
$ws.Range("N65").Value = -78737.5

$ws.Range("H136").Value = 473.3
$ws.Range("I136").Value = 447.875
$ws.Range("J136").Value = 575
$ws.Range("K136").Value = 1343.625
$ws.Range("L136").Value = 1725
$ws.Range("M136").Value = 1206.375
$ws.Range("N136").Value = -6825

Write-Output "Applied all cell updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR"
